# New Org Production Debugging
# - product_catalog_test_data: "Alexa67 : 2" -> "Alexa678 : 2" (buy_items, col M)
# - product_catalog_test_data: "Ranger Alex" -> "Alexa678"     (product_name/price_book_name)
# - Sheet1 window selection moved to AC18 (scrolled right to show column U onward)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update every cell that shared the "Alexa67 : 2" string ---
$ws.Range("M7").Value  = "Alexa678 : 2"
$ws.Range("M9").Value  = "Alexa678 : 2"
$ws.Range("M10").Value = "Alexa678 : 2"
$ws.Range("M11").Value = "Alexa678 : 2"

# --- Update every cell that shared the "Ranger Alex" string ---
$ws.Range("AF8").Value  = "Alexa678"
$ws.Range("AF9").Value  = "Alexa678"
$ws.Range("V10").Value  = "Alexa678"
$ws.Range("AF10").Value = "Alexa678"
$ws.Range("AF11").Value = "Alexa678"

# --- Move the sheet's selection/active cell to AC18 (as last edited by the author) ---
$ws.Range("AC18").Select()
